$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.102666333333334
$ws.Range("H2").Value = 6.307999000000001
$ws.Range("I2").Value = 0.03992401473981187
$ws.Range("J2").Value = 0.03992401473981187
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.067615666666667
$ws.Range("N2").Value = 15.202847
$ws.Range("O2").Value = 0.6207828410514926
$ws.Range("P2").Value = 0.6207828410514925
$ws.Range("Q2").Value = 10.65550485257256
$ws.Range("R2").Value = 95.899543673153
$ws.Range("S2").Value = 0.02478414329636208
$ws.Range("T2").Value = 0.02478414329636208
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.102666333333334
$ws.Range("H3").Value = 6.307999000000001
$ws.Range("I3").Value = 0.03992401473981187
$ws.Range("J3").Value = 0.03992401473981187
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.388310333333333
$ws.Range("N3").Value = 4.164931
$ws.Range("O3").Value = 0.1700679944331107
$ws.Range("P3").Value = 0.1700679944331107
$ws.Range("Q3").Value = 2.919153398118778
$ws.Range("R3").Value = 26.272380583069
$ws.Range("S3").Value = 0.006789797116517756
$ws.Range("T3").Value = 0.006789797116517755
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.102666333333334
$ws.Range("H4").Value = 6.307999000000001
$ws.Range("I4").Value = 0.03992401473981187
$ws.Range("J4").Value = 0.03992401473981187
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.707340333333333
$ws.Range("N4").Value = 5.122021
$ws.Range("O4").Value = 0.2091491645153968
$ws.Range("P4").Value = 0.2091491645153968
$ws.Range("Q4").Value = 3.589967038442112
$ws.Range("R4").Value = 32.309703345979
$ws.Range("S4").Value = 0.008350074326932041
$ws.Range("T4").Value = 0.008350074326932039
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.75404733333334
$ws.Range("H5").Value = 83.26214200000001
$ws.Range("I5").Value = 0.526975192050016
$ws.Range("J5").Value = 0.5269751920500161
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.067615666666667
$ws.Range("N5").Value = 15.202847
$ws.Range("O5").Value = 0.6207828410514926
$ws.Range("P5").Value = 0.6207828410514925
$ws.Range("Q5").Value = 140.6468450798082
$ws.Range("R5").Value = 1265.821605718274
$ws.Range("S5").Value = 0.3271371568844649
$ws.Range("T5").Value = 0.3271371568844649
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 27.75404733333334
$ws.Range("H6").Value = 83.26214200000001
$ws.Range("I6").Value = 0.526975192050016
$ws.Range("J6").Value = 0.5269751920500161
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.388310333333333
$ws.Range("N6").Value = 4.164931
$ws.Range("O6").Value = 0.1700679944331107
$ws.Range("P6").Value = 0.1700679944331107
$ws.Range("Q6").Value = 38.53123070468911
$ws.Range("R6").Value = 346.7810763422021
$ws.Range("S6").Value = 0.08962161402794956
$ws.Range("T6").Value = 0.08962161402794958
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 27.75404733333334
$ws.Range("H7").Value = 83.26214200000001
$ws.Range("I7").Value = 0.526975192050016
$ws.Range("J7").Value = 0.5269751920500161
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.707340333333333
$ws.Range("N7").Value = 5.122021
$ws.Range("O7").Value = 0.2091491645153968
$ws.Range("P7").Value = 0.2091491645153968
$ws.Range("Q7").Value = 47.38560442544245
$ws.Range("R7").Value = 426.4704398289821
$ws.Range("S7").Value = 0.1102164211376016
$ws.Range("T7").Value = 0.1102164211376016
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 22.809992
$ws.Range("H8").Value = 68.42997600000001
$ws.Range("I8").Value = 0.4331007932101721
$ws.Range("J8").Value = 0.4331007932101721
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.067615666666667
$ws.Range("N8").Value = 15.202847
$ws.Range("O8").Value = 0.6207828410514926
$ws.Range("P8").Value = 0.6207828410514925
$ws.Range("Q8").Value = 115.5922728157414
$ws.Range("R8").Value = 1040.330455341672
$ws.Range("S8").Value = 0.2688615408706657
$ws.Range("T8").Value = 0.2688615408706656
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 22.809992
$ws.Range("H9").Value = 68.42997600000001
$ws.Range("I9").Value = 0.4331007932101721
$ws.Range("J9").Value = 0.4331007932101721
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.388310333333333
$ws.Range("N9").Value = 4.164931
$ws.Range("O9").Value = 0.1700679944331107
$ws.Range("P9").Value = 0.1700679944331107
$ws.Range("Q9").Value = 31.66734759685067
$ws.Range("R9").Value = 285.0061283716561
$ws.Range("S9").Value = 0.07365658328864338
$ws.Range("T9").Value = 0.07365658328864338
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 22.809992
$ws.Range("H10").Value = 68.42997600000001
$ws.Range("I10").Value = 0.4331007932101721
$ws.Range("J10").Value = 0.4331007932101721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.707340333333333
$ws.Range("N10").Value = 5.122021
$ws.Range("O10").Value = 0.2091491645153968
$ws.Range("P10").Value = 0.2091491645153968
$ws.Range("Q10").Value = 38.94441934461068
$ws.Range("R10").Value = 350.4997741014961
$ws.Range("S10").Value = 0.09058266905086315
$ws.Range("T10").Value = 0.09058266905086314
